$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.81"
$ws.Range("E2").Value = "'7.08%"
$ws.Range("G2").Value = "'11"
$ws.Range("D3").Value = "'48.75"
$ws.Range("E3").Value = "'15.48%"
$ws.Range("G3").Value = "'11"
$ws.Range("D4").Value = "'5.277"
$ws.Range("E4").Value = "'5.24%"
$ws.Range("G4").Value = "'11"
$ws.Range("D5").Value = "'0.08110"
$ws.Range("E5").Value = "'7.30%"
$ws.Range("G5").Value = "'11"
$ws.Range("D6").Value = "'4.600"
$ws.Range("E6").Value = "'5.07%"
$ws.Range("G6").Value = "'11"
$ws.Range("D7").Value = "'1.646"
$ws.Range("E7").Value = "'2.68%"
$ws.Range("G7").Value = "'11"
$ws.Range("E8").Value = "'28.87%"
$ws.Range("G8").Value = "'11"
$ws.Range("D9").Value = "'0.1293"
$ws.Range("E9").Value = "'8.59%"
$ws.Range("G9").Value = "'11"
$ws.Range("D10").Value = "'0.1939"
$ws.Range("E10").Value = "'5.26%"
$ws.Range("G10").Value = "'11"
$ws.Range("D11").Value = "'0.09439"
$ws.Range("E11").Value = "'3.94%"
$ws.Range("G11").Value = "'11"
$ws.Range("D12").Value = "'0.04589"
$ws.Range("E12").Value = "'9.97%"
$ws.Range("G12").Value = "'11"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("E13").Value = "'0.22%"
$ws.Range("G13").Value = "'11"
$ws.Range("D14").Value = "'0.001325"
$ws.Range("E14").Value = "'3.29%"
$ws.Range("G14").Value = "'11"
$ws.Range("D15").Value = "'0.04171"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("G15").Value = "'11"
$ws.Range("D16").Value = "'0.005851"
$ws.Range("E16").Value = "'-0.88%"
$ws.Range("G16").Value = "'11"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("G17").Value = "'11"
$ws.Range("D18").Value = "'2.431"
$ws.Range("E18").Value = "'1.94%"
$ws.Range("G18").Value = "'11"
$ws.Range("D19").Value = "'0.3422"
$ws.Range("E19").Value = "'2.61%"
$ws.Range("G19").Value = "'11"
$ws.Range("D20").Value = "'8.123"
$ws.Range("E20").Value = "'-3.31%"
$ws.Range("G20").Value = "'11"
$ws.Range("D21").Value = "'0.1391"
$ws.Range("E21").Value = "'-1.25%"
$ws.Range("G21").Value = "'11"
$ws.Range("E22").Value = "'-5.21%"
$ws.Range("G22").Value = "'11"
$ws.Range("D23").Value = "'0.001305"
$ws.Range("E23").Value = "'3.27%"
$ws.Range("G23").Value = "'11"
$ws.Range("D24").Value = "'0.004252"
$ws.Range("E24").Value = "'8.88%"
$ws.Range("G24").Value = "'11"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("E25").Value = "'6.53%"
$ws.Range("G25").Value = "'11"
$ws.Range("D26").Value = "'0.0003541"
$ws.Range("G26").Value = "'11"
$ws.Range("G27").Value = "'11"
$ws.Range("G28").Value = "'11"
$ws.Range("G29").Value = "'11"
$ws.Range("G30").Value = "'11"
$ws.Range("G31").Value = "'11"
$ws.Range("G32").Value = "'11"
$ws.Range("G33").Value = "'11"
$ws.Range("G34").Value = "'11"
$ws.Range("G35").Value = "'11"
$ws.Range("G36").Value = "'11"
$ws.Range("G37").Value = "'11"
$ws.Range("D38").Value = "'0.02718"
$ws.Range("E38").Value = "'12.44%"
$ws.Range("G38").Value = "'11"
$ws.Range("D39").Value = "'0.05701"
$ws.Range("E39").Value = "'9.10%"
$ws.Range("G39").Value = "'11"
$ws.Range("D40").Value = "'0.006302"
$ws.Range("E40").Value = "'-7.22%"
$ws.Range("G40").Value = "'11"
$ws.Range("D41").Value = "'0.007705"
$ws.Range("E41").Value = "'0.06%"
$ws.Range("G41").Value = "'11"
$ws.Range("D42").Value = "'0.1441"
$ws.Range("E42").Value = "'8.29%"
$ws.Range("G42").Value = "'11"
$ws.Range("D43").Value = "'0.007681"
$ws.Range("E43").Value = "'4.08%"
$ws.Range("G43").Value = "'11"
$ws.Range("E44").Value = "'3.93%"
$ws.Range("G44").Value = "'11"
$ws.Range("E45").Value = "'6.41%"
$ws.Range("G45").Value = "'11"
$ws.Range("D46").Value = "'0.00006928"
$ws.Range("E46").Value = "'10.92%"
$ws.Range("G46").Value = "'11"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("G47").Value = "'11"
$ws.Range("D48").Value = "'0.05502"
$ws.Range("E48").Value = "'20.39%"
$ws.Range("G48").Value = "'11"
$ws.Range("D49").Value = "'0.004002"
$ws.Range("G49").Value = "'11"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("G50").Value = "'11"
$ws.Range("E51").Value = "'0.11%"
$ws.Range("G51").Value = "'11"
